$d = $word.ActiveDocument

# Locate the "Implement full set" paragraph (rather than hard-coding an
# index) so the script is resilient to minor structural differences.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Implement full set") {
        $anchorIndex = $i
    }
}

# --- New paragraph: "Implement mortgage" ---------------------------------
# Typed as two runs ("Implement m" + "ortgage"), matching how Word split
# the text across two separate insertions/edits.
$firstPart = "Implement m"
$secondPart = "ortgage"

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter()

$mortgageIndex = $anchorIndex + 1
$pMortgage = $d.Paragraphs.Item($mortgageIndex)
$pMortgage.Range.InsertAfter($firstPart)

$pMortgage2 = $d.Paragraphs.Item($mortgageIndex)
$pMortgage2.Range.InsertAfter($secondPart)

# Toggle Bold on/off over the just-appended "ortgage" text. This forces the
# engine to keep it as a distinct run object even though the final
# formatting (Bold off) matches the preceding run exactly.
$pMortgage3 = $d.Paragraphs.Item($mortgageIndex)
$paraEnd = $pMortgage3.Range.End - 1
$boundary = $paraEnd - $secondPart.Length
$appended = $d.Range($boundary, $paraEnd)
$appended.Bold = 1
$appended.Bold = 0

# --- New paragraph: "Implement trading" -----------------------------------
$pMortgage4 = $d.Paragraphs.Item($mortgageIndex)
$pMortgage4.Range.InsertParagraphAfter()

$tradingIndex = $mortgageIndex + 1
$pTrading = $d.Paragraphs.Item($tradingIndex)
$pTrading.Range.InsertAfter("Implement trading")

# --- New paragraph: "Make houses and hotels" ------------------------------
$pTrading2 = $d.Paragraphs.Item($tradingIndex)
$pTrading2.Range.InsertParagraphAfter()

$housesIndex = $tradingIndex + 1
$pHouses = $d.Paragraphs.Item($housesIndex)
$pHouses.Range.InsertAfter("Make houses and hotels")

Write-Output "done"
